# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for the latest scrape, per the GitHub Actions cron job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.800.58'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '3.161.42'
$ws.Range("E3").Value = '  +2.28%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.09'
$ws.Range("E5").Value = '  +2.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.66'
$ws.Range("E6").Value = '  +5.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.157.61'
$ws.Range("E8").Value = '  +2.33%  '
$ws.Range("E9").Value = '  +4.84%  '
$ws.Range("E10").Value = '  +6.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.16'
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("E12").Value = '  +7.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000258'
$ws.Range("E13").Value = '  +13.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.15'
$ws.Range("E14").Value = '  +8.78%  '
$ws.Range("D15").Value = '3.677.00'
$ws.Range("E15").Value = '  +2.52%  '
$ws.Range("D16").Value = '64.904.72'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.27'
$ws.Range("E17").Value = '  +7.74%  '
$ws.Range("D18").Value = '3.160.96'
$ws.Range("E18").Value = '  +2.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '520.98'
$ws.Range("E19").Value = '  +7.29%  '
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.02'
$ws.Range("E21").Value = '  +7.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.740'
$ws.Range("E22").Value = '  +9.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.27'
$ws.Range("E23").Value = '  +8.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.87'
$ws.Range("E24").Value = '  +3.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.45'
$ws.Range("E25").Value = '  +5.21%  '
$ws.Range("E27").Value = '  +5.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.80'
$ws.Range("E28").Value = '  +10.09%  '
$ws.Range("E29").Value = '  +6.96%  '
$ws.Range("E30").Value = '  +6.38%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.68'
$ws.Range("E32").Value = '  +7.99%  '
$ws.Range("E33").Value = '  +3.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.16'
$ws.Range("E34").Value = '  +10.35%  '
$ws.Range("E35").Value = '  +6.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.81'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '489.66'
$ws.Range("E37").Value = '  +8.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0869'
$ws.Range("E38").Value = '  +6.34%  '
$ws.Range("E39").Value = '  +4.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.00'
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("D41").Value = '3.115.73'
$ws.Range("E41").Value = '  +4.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.69'
$ws.Range("E42").Value = '  +5.49%  '
$ws.Range("E43").Value = '  +14.61%  '
$ws.Range("E44").Value = '  +6.71%  '
$ws.Range("E45").Value = '  +16.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.37'
$ws.Range("E46").Value = '  +5.24%  '
$ws.Range("E47").Value = '  +13.26%  '
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("E49").Value = '  +3.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.31'
$ws.Range("E50").Value = '  +11.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.44'
$ws.Range("E51").Value = '  +0.28%  '
